# Replace each arithmetic-answer cell of the single 20x5 table with its
# updated equation, in row-major (left-to-right, top-to-bottom) order,
# matching the order the cells already appear in the document.
$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$rows = $t.Rows.Count
$cols = $t.Columns.Count

$newValues = @(
  "26+46=72",
  "39+17=56",
  "85+9=94",
  "19+19=38",
  "3+9=12",
  "50-12=38",
  "97-19=78",
  "55-6=49",
  "61-39=22",
  "29+47=76",
  "66+16=82",
  "92-76=16",
  "38+44=82",
  "73-29=44",
  "91-58=33",
  "40-22=18",
  "93-78=15",
  "85-37=48",
  "59+4=63",
  "9+62=71",
  "71-42=29",
  "57+24=81",
  "19+79=98",
  "19+25=44",
  "4+39=43",
  "78+13=91",
  "48+43=91",
  "45+39=84",
  "17+18=35",
  "16+75=91",
  "61-53=8",
  "63-24=39",
  "19+27=46",
  "43-36=7",
  "81-2=79",
  "36+7=43",
  "47+37=84",
  "26+6=32",
  "77+8=85",
  "26+29=55",
  "94-76=18",
  "9+76=85",
  "72-65=7",
  "33-29=4",
  "83-5=78",
  "84-59=25",
  "37-29=8",
  "73-38=35",
  "95-78=17",
  "74-48=26",
  "46+37=83",
  "90-32=58",
  "74-25=49",
  "6+49=55",
  "48+23=71",
  "61-6=55",
  "69+4=73",
  "61-44=17",
  "14-8=6",
  "82-17=65",
  "30-26=4",
  "9+68=77",
  "62-16=46",
  "9+57=66",
  "25+47=72",
  "19+55=74",
  "58+19=77",
  "70-41=29",
  "86-58=28",
  "59+22=81",
  "81-27=54",
  "36+28=64",
  "45+37=82",
  "47+6=53",
  "90-66=24",
  "33-9=24",
  "41-3=38",
  "57+28=85",
  "51-22=29",
  "16+29=45",
  "68+28=96",
  "23+8=31",
  "88+5=93",
  "80-2=78",
  "11-2=9",
  "64-47=17",
  "16+27=43",
  "70-61=9",
  "84-16=68",
  "87+7=94",
  "51-44=7",
  "38+18=56",
  "40-11=29",
  "61-22=39",
  "71-59=12",
  "18+19=37",
  "37+5=42",
  "30-7=23",
  "78+13=91",
  "51-36=15"
)

$idx = 0
for ($r = 1; $r -le $rows; $r++) {
  for ($c = 1; $c -le $cols; $c++) {
    $cell = $t.Cell($r, $c)
    $cell.Range.Text = $newValues[$idx]
    $idx = $idx + 1
  }
}